# ActivityLogSheetWk9 - "moved local copies to repo"
# Fill in week-9 activity log entries (name, week number, activities,
# dates/times, group vs individual hours) that were populated once the
# workbook moved from the author's local machine into the repo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block -----------------------------------------------------
$ws.Range("B2").Value = "Jesse Hare"
$ws.Range("G2").Value = 9

# Dates/times are written as raw Excel serial numbers (days since
# 1899-12-30, fractional = time-of-day) so they land as real numeric
# cells using the date/time number formats already on these cells,
# rather than as text.

# --- Activity rows (4-9) ------------------------------------------------
# Row 4: Review progress - Group - 23/09/2019 09:00 -> 10:00, 1 group hr
$ws.Range("A4").Value = "Review progress"
$ws.Range("C4").Value = "G"
$ws.Range("D4").Value = 43731
$ws.Range("E4").Value = 0.375
$ws.Range("F4").Value = 0.41666666666666669
$ws.Range("G4").Value = 1

# Row 5: Prepare for next meeting - Group - 24/09/2019 09:00 -> 11:00, 2 group hrs
$ws.Range("A5").Value = "Prepare for next meeting"
$ws.Range("C5").Value = "G"
$ws.Range("D5").Value = 43732
$ws.Range("E5").Value = 0.375
$ws.Range("F5").Value = 0.45833333333333331
$ws.Range("G5").Value = 2

# Row 6: Plan next iteration - Group - 25/09/2019 09:00 -> 11:00, 2 group hrs
$ws.Range("A6").Value = "Plan next iteration"
$ws.Range("C6").Value = "G"
$ws.Range("D6").Value = 43733
$ws.Range("E6").Value = 0.375
$ws.Range("F6").Value = 0.45833333333333331
$ws.Range("G6").Value = 2

# Row 7: Work on next iteration - Individual - 26/09/2019 09:00 -> 02:00, 5 individual hrs
$ws.Range("A7").Value = "Work on next iteration"
$ws.Range("C7").Value = "I"
$ws.Range("D7").Value = 43734
$ws.Range("E7").Value = 0.375
$ws.Range("F7").Value = 0.083333333333333329
$ws.Range("H7").Value = 5

# Row 8: Work on next iteration - Individual - 27/09/2019 09:00 -> 15:00, 6 individual hrs
$ws.Range("A8").Value = "Work on next iteration"
$ws.Range("C8").Value = "I"
$ws.Range("D8").Value = 43735
$ws.Range("E8").Value = 0.375
$ws.Range("F8").Value = 0.625
$ws.Range("H8").Value = 6

# Row 9: Work on next iteration - Individual - 28/09/2019 09:00 -> 01:00, 4 individual hrs
$ws.Range("A9").Value = "Work on next iteration"
$ws.Range("C9").Value = "I"
$ws.Range("D9").Value = 43736
$ws.Range("E9").Value = 0.375
$ws.Range("F9").Value = 0.041666666666666664
$ws.Range("H9").Value = 4

# --- Column widths -------------------------------------------------------
# Column B (Name) widened to fit "Jesse Hare"/longer text, D (Date) and F
# (End Time) widened to best-fit the date/time values just entered.
$ws.Columns("B").ColumnWidth = 24.75
$ws.Columns("D").ColumnWidth = 12.1875
$ws.Columns("F").ColumnWidth = 12.75

# --- Selection moved to A5:B5 as a final cursor position ---------------
$ws.Range("A5:B5").Select()
